$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# C10 was 18, update it to 1 (numeric value).
$ws.Range("C10").Value = 1
